$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: Iterations = 100 (rows 2-11) ---

# Row 2 (Nelson24)
$ws.Range("D2").Value = 0.0001205196183317144
$ws.Range("E2").Value = 0.0001205196183317144

# Row 3 (Nelson30)
$ws.Range("D3").Value = [double]"4.41755591186678E-07"
$ws.Range("E3").Value = [double]"4.41755591186678E-07"

# Row 4 (Nelson57)
$ws.Range("D4").Value = 0.001863934658813622
$ws.Range("E4").Value = 0.001863934658813622

# Row 5 (Nelson49)
$ws.Range("D5").Value = 0.003822217619246418
$ws.Range("E5").Value = 0.003822217619246418

# Row 6 (Nelson53)
$ws.Range("D6").Value = 0.1654206248896664
$ws.Range("E6").Value = 0.1654206248896664

# Row 7 (Ableson 0)
$ws.Range("D7").Value = 0.9925850874805752
$ws.Range("E7").Value = 0.007414912519424832

# Row 8 (Ableson 12)
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.08589535494779904
$ws.Range("E8").Value = 0.914104645052201

# Row 9 (Ableson 11)
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.0006806875824239165
$ws.Range("E9").Value = 0.9993193124175761

# Row 10 (Ableson 17)
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.05198062877698079
$ws.Range("E10").Value = 0.9480193712230192

# Row 11 (Ableson 44)
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.01026406553360165
$ws.Range("E11").Value = 0.9897359344663983
$ws.Range("F11").Value = 1.747710943222046
$ws.Range("G11").Value = 0.6

# --- Block 2: Iterations = 200 (rows 12-21) ---

# Row 12 (Nelson24)
$ws.Range("D12").Value = [double]"7.333622789178556E-07"
$ws.Range("E12").Value = [double]"7.333622789178556E-07"

# Row 13 (Nelson30)
$ws.Range("D13").Value = [double]"2.603280256073345E-09"
$ws.Range("E13").Value = [double]"2.603280256073345E-09"

# Row 14 (Nelson57)
$ws.Range("D14").Value = 0.0002217453326150737
$ws.Range("E14").Value = 0.0002217453326150737

# Row 15 (Nelson49)
$ws.Range("D15").Value = 0.0002114032010213806
$ws.Range("E15").Value = 0.0002114032010213806

# Row 16 (Nelson53)
$ws.Range("D16").Value = 0.03925039207597144
$ws.Range("E16").Value = 0.03925039207597144

# Row 17 (Ableson 0)
$ws.Range("D17").Value = 0.9588862564789086
$ws.Range("E17").Value = 0.04111374352109143

# Row 18 (Ableson 12)
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.06508673683908942
$ws.Range("E18").Value = 0.9349132631609106

# Row 19 (Ableson 11)
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.0006413738786290194
$ws.Range("E19").Value = 0.9993586261213709

# Row 20 (Ableson 17)
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.03014554556788774
$ws.Range("E20").Value = 0.9698544544321123

# Row 21 (Ableson 44)
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.002691205399276805
$ws.Range("E21").Value = 0.9973087946007232
$ws.Range("F21").Value = 1.958587408065796
$ws.Range("G21").Value = 0.6
